# "adding profits to tables"
# Adds a third block of columns (R:Y) mirroring the existing GFA/IMF/OECD
# sub-header layout already present in B:I and J:Q, but labelled "M_PL"
# (a new top header value) and populated with the new profit figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Top header (row 1): merged "M_PL" label over R1:Y1, mirroring the
#     existing "M_%cit" (B1:I1) / "M_ETR" (J1:Q1) merged headers.
$ws.Range("R1").Value = "M_PL"
$ws.Range("R1:Y1").Merge()

$headerRange = $ws.Range("R1:Y1")
$headerRange.Borders.LineStyle = 1
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop

# --- Sub-header (row 2): same 8 sub-column labels repeated under the new
#     "M_PL" block, exactly as they already repeat under M_%cit and M_ETR.
$subHeaders = @("GFA - Sales", "GFA - Sales + Emp", "IMF - Sales", "IMF - Sales + Emp", "OECD (20%) - Sales", "OECD (20%) - Sales + Emp", "OECD - Sales", "OECD - Sales + Emp")
$subCols = @("R", "S", "T", "U", "V", "W", "X", "Y")

for ($i = 0; $i -lt $subCols.Length; $i++) {
    $cell = $ws.Range($subCols[$i] + "2")
    $cell.Value = $subHeaders[$i]
    $cell.Borders.LineStyle = 1
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

# --- Data rows 4-8: new profit-level ("M_PL") figures for each of the 8
#     sub-columns, one row per group (High Income / LICs / LMICs / Tax
#     haven / UMICs).
$data = @{
    4 = @(957691691302, 957942582918, 955500195836, 955751087452, 1007534436142, 1007534436142, 1007534436142, 1007534436142)
    5 = @(-269766813, -269766813, -269766813, -269766813, -269766813, -269766813, -269766813, -269766813)
    6 = @(1507896554, 15012781826, 1074779159, 16007585528, 20228669958, 20228669958, 20228669958, 20228669958)
    7 = @(2770915272, -7824629507, -6943500091, -7824629507, -8693540732, -8693540732, -8693540732, -8693540732)
    8 = @(37772977295, 43238497013, 37492196667, 43740019703, 45733381438, 45733381438, 45733381438, 45733381438)
}

foreach ($row in $data.Keys) {
    $rowValues = $data[$row]
    for ($i = 0; $i -lt $subCols.Length; $i++) {
        $ws.Range($subCols[$i] + $row).Value = $rowValues[$i]
    }
}
